$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Objective changes from "Defensive" to "Standard" (still a "View" row for Price to Book)
$ws.Range("B8").Value = "Standard"

# Row 9: was a duplicate "View / Standard / Price to Book" row; replace it with a
# "Search / Standard / Market Cap" row (mirrors the old row 11 content).
$ws.Range("A9").Value = "Search"
$ws.Range("B9").Value = "Standard"
$ws.Range("C9").Value = "Market Cap"
$ws.Range("D9").Value = "marketcap"
$ws.Range("E9").Value = "~gt~"
$ws.Range("F9").Value = 2000000000

# Row 10: was "View / Risky / Price to Revenue"; replace it with
# "View / Standard / Free Cash Flow to Firm" (mirrors the old row 12 content).
$ws.Range("A10").Value = "View"
$ws.Range("B10").Value = "Standard"
$ws.Range("C10").Value = "Free Cash Flow to Firm"
$ws.Range("D10").Value = "fcffgrowth"
$ws.Range("E10").Value = "~gt~"
$ws.Range("F10").Value = -99

# Rows 11-13 are no longer part of the table; delete them entirely.
$ws.Range("A13:H13").EntireRow.Delete()
$ws.Range("A12:H12").EntireRow.Delete()
$ws.Range("A11:H11").EntireRow.Delete()

# Move the active selection to reflect the new, shorter table (next empty row area).
$ws.Range("E14").Select()
